# Add rows via upload: update existing rows 5-8 and append new rows 9-18
# to the data table on Sheet1, extending the used range from A1:K8 to A1:K18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, followed by column A..K values for that row.
# H column sometimes holds the literal text "inf" (urea/NaOH ratio when
# NaOH volume is 0) instead of a number.
$rows = @(
    @(5,  3,  1.2, 0.2, 0.2, 44, 0.9, 0.9, 1,     1,   2.2, 85),
    @(6,  6,  1.2, 0.2, 0.2, 43, 0.9, 1.9, 0.5,   1.2, 2.2, 82),
    @(7,  7,  1.2, 0.3, 0.1, 44, 0.6, 1.4, 0.4,   1,   2.2, 72),
    @(8,  9,  1.2, 0.2, 0.2, 44, 1.3, 1.4, 0.9,   1,   2.2, 74),
    @(9,  10, 1.2, 0.2, 0.2, 45, 0.2, 1,   0.2,   1.2, 2.2, 83),
    @(10, 11, 1.2, 0.1, 0.3, 44, 0.8, 1.4, 0.6,   1.2, 2.2, 84),
    @(11, 12, 1.2, 0.2, 0.2, 45, 1.1, 0,   "inf", 1.1, 2.2, 84),
    @(12, 14, 1.2, 0.3, 0.1, 43, 1.1, 1.8, 0.6,   0.9, 2.2, 84),
    @(13, 15, 1.2, 0.1, 0.3, 44, 1.4, 0.7, 2,     1.2, 2.2, 84),
    @(14, 17, 1.2, 0.2, 0.2, 44, 0,   1.7, 0,     1.3, 2.2, 66),
    @(15, 19, 1.2, 0.1, 0.3, 46, 0.4, 0,   "inf", 1.2, 2.2, 85),
    @(16, 20, 1.2, 0.2, 0.2, 46, 0.1, 0,   "inf", 1,   2.2, 83),
    @(17, 24, 1.2, 0.2, 0.2, 45, 1.7, 0,   "inf", 1,   2.2, 84),
    @(18, 27, 1.2, 0.2, 0.2, 42, 0.6, 2.8, 0.2,   1.3, 2.2, 84)
)

foreach ($row in $rows) {
    $r = $row[0]
    for ($col = 1; $col -le 11; $col++) {
        $ws.Cells.Item($r, $col).Value = $row[$col]
    }
}
